# Update countries & provincias Spain
#
# The underlying data source re-sorted rows that are tied on total
# cases (column B), which changes the relative order of a handful of
# countries that had identical stats. Column A (country name) and the
# B:H statistic columns always travel together as one row, so each
# affected destination row is rewritten with the source row's full
# tuple of values.
#
# Set-CountryRow Row Name B C D E F G H

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CountryRow($Row, $Name, $B, $C, $D, $E, $F, $G, $H) {
    $ws.Cells.Item($Row, 1).Value = $Name
    $ws.Cells.Item($Row, 2).Value = $B
    $ws.Cells.Item($Row, 3).Value = $C
    $ws.Cells.Item($Row, 4).Value = $D
    $ws.Cells.Item($Row, 5).Value = $E
    $ws.Cells.Item($Row, 6).Value = $F
    $ws.Cells.Item($Row, 7).Value = $G
    $ws.Cells.Item($Row, 8).Value = $H
}

# Cycle: Bahamas / Guinea Ecuatorial / Birmania
Set-CountryRow 156 "Guinea Ecuatorial" 15 0 1 14 0 0 0
Set-CountryRow 157 "Birmania"          15 0 0 14 0 0 1
Set-CountryRow 159 "Bahamas"           15 1 1 14 0 0 0

# Cycle: Mongolia / Guyana
Set-CountryRow 164 "Guyana"            12 0 0 10 0 0 2
Set-CountryRow 165 "Mongolia"          12 0 2 10 0 0 0

# Cycle: Siria / Groenlandia
Set-CountryRow 171 "Groenlandia"       10 0 2 8  0 0 0
Set-CountryRow 172 "Siria"             10 0 0 8  0 0 2

# Cycle: Suazilandia / Granada
Set-CountryRow 174 "Granada"           9  0 0 9  0 0 0
Set-CountryRow 175 "Suazilandia"       9  0 0 9  0 0 0

# Cycle: Angola / Sudan
Set-CountryRow 183 "Sudan"             7  0 1 4  0 0 2
Set-CountryRow 184 "Angola"            7  0 1 4  0 0 2

# Cycle: Santa Sede / San Martin (Parte Holandesa)
Set-CountryRow 185 "San Martin (Parte Holandesa)" 6 0 0 6 0 0 0
Set-CountryRow 186 "Santa Sede"         6  0 0 6  0 0 0

# Cycle: San Bartolome / Cabo Verde
Set-CountryRow 187 "Cabo Verde"         6  0 0 5  0 0 1
Set-CountryRow 188 "San Bartolome"      6  0 1 5  0 0 0

# Cycle: Fiyi / Montserrat / Islas Turcas y Caicos
Set-CountryRow 190 "Montserrat"             5 0 0 5 0 0 0
Set-CountryRow 191 "Islas Turcas y Caicos"  5 0 0 5 0 0 0
Set-CountryRow 192 "Fiyi"                   5 0 0 5 0 0 0

# Cycle: Nepal / Somalia
Set-CountryRow 193 "Somalia"           5  0 1 4  0 0 0
Set-CountryRow 195 "Nepal"             5  0 1 4  0 0 0

# Cycle: Belice / Republica de Africa Central / Liberia / Islas Virgenes Britanicas
Set-CountryRow 199 "Republica de Africa Central" 3 0 0 3 0 0 0
Set-CountryRow 200 "Liberia"                     3 0 0 3 0 0 0
Set-CountryRow 201 "Islas Virgenes Britanicas"   3 0 0 3 0 0 0
Set-CountryRow 202 "Belice"                      3 0 0 3 0 0 0

# Cycle: Sierra Leona / Papua Nueva Guinea
Set-CountryRow 207 "Papua Nueva Guinea" 1  0 0 1  0 0 0
Set-CountryRow 208 "Sierra Leona"       1  0 0 1  0 0 0

# Update the "last refreshed" timestamp banner in A1.
$ws.Range("A1").Value = "Datos actualizados a 1 de Abril de 2020 a las 06:20"
